# Applies the content update described by the commit:
# "Created a file to sort all of the papers into respective folders.
#  Changed GPT response from list to delimited string"
#
# Row 1 = data for the "66266 optocoupler" paper (Micropac Industries, TID test)
# Row 2 = data for the "SEE hardened power management ICs" paper (Intersil, SEE test)
# Row 2's trailing details (radiation type / energy / fluence / failures), which used
# to be split across columns F/G/H/I, are now a single delimited string in column F,
# with the SEE/TID classification moved to column G - H2/I2 are no longer used.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear any stale content first (columns A:I, rows 1:2) and rewrite from scratch.
$ws.Range("A1:I2").ClearContents()

# Cell C1 ("66266") looks numeric; force it to be treated as text first so that
# the leading/trailing spaces and the textual representation are preserved
# instead of Excel auto-converting it into the number 66266.
$ws.Cells.Item(1, 3).NumberFormat = "@"

# Values are written in this particular order (file names first, then the rest
# of row 1, then the rest of row 2) so that the shared-string table ends up
# populated in the same sequence as the source data.
$ws.Cells.Item(1, 1).Value = "Papers_Sorted/SMD\Characterization_of_the_Effects_of_250_MeV_Proton-Induced_Total_Ionizing_Dose_and_Displacement_Damage_on_the_66266_Optocoupler.pdf"
$ws.Cells.Item(2, 1).Value = "Papers_Sorted/SMD\Characterization_of_various_SEE_hardened_power_management_ICs.pdf"

# Row 1
$ws.Cells.Item(1, 2).Value = "S. Messenger "
$ws.Cells.Item(1, 3).Value = " 66266 "
$ws.Cells.Item(1, 4).Value = " optocoupler "
$ws.Cells.Item(1, 5).Value = " Micropac Industries, Inc. "
$ws.Cells.Item(1, 6).Value = " TID"
$ws.Cells.Item(1, 7).Value = "Type of radiation source: 250 MeV protons "
$ws.Cells.Item(1, 8).Value = " Failures: None "
$ws.Cells.Item(1, 9).Value = " When: N/A"

# Drop the auxiliary text-format style so the cell keeps no explicit style index.
$ws.Cells.Item(1, 3).Style = "Normal"

# Row 2
$ws.Cells.Item(2, 2).Value = "B. P. Alaskiewicz "
$ws.Cells.Item(2, 3).Value = " IS-2100ARH "
$ws.Cells.Item(2, 4).Value = " MOSFET driver IC "
$ws.Cells.Item(2, 5).Value = " Intersil Corporation "
$ws.Cells.Item(2, 6).Value = " SEE"
$ws.Cells.Item(2, 7).Value = "Radiation Source Type: Au ions at 90.9MeV/mg/cm²`nFailures: No`n"

# The embedded line breaks above make Excel mark row 2 with an explicit custom
# row height; auto-fit it back down so no extra height metadata is persisted.
$ws.Rows.Item(2).AutoFit()
